# feat: save employee branch in import employee
#
# Adds a new "Branch" column (W) to the employee-list import template,
# matching the header formatting of the neighbouring "Works days per
# year" column (V), widens the new column, updates the date format
# used by the "Birth date" / "Date employed" columns, and restores the
# sheet's scroll/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Branch" header cell in column W -------------------
$ws.Range("W1").Value = "Branch"

# Match the bold header formatting used by the rest of row 1 (e.g. V1)
$ws.Range("W1").Font.Name = "Calibri"
$ws.Range("W1").Font.Size = 11
$ws.Range("W1").Font.Bold = $true
$ws.Range("W1").NumberFormat = "General"
$ws.Range("W1").Locked = $true

# --- Give the new column its own width -------------------------------
$ws.Columns.Item(23).ColumnWidth = 14.6

# --- Update the date format used for Birth date / Date employed ------
$ws.Columns.Item(6).NumberFormat = "M/D/YYYY"
$ws.Columns.Item(18).NumberFormat = "M/D/YYYY"

# --- Restore view / selection state -----------------------------------
[void]$ws.Activate()
[void]$ws.Range("S1").Select()
[void]$ws.Range("W7").Select()
